$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-120 down to 27-121.
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new data record.
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 45251
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 100112012
$ws.Range("G26").Value = "Espinaca"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 550
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 945
$ws.Range("N26").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 315
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = "Hortaliza"
